# Weekly update: insert a new price record for Espinaca (Vega Central Mapocho
# de Santiago) at row 276, pushing all following rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 276 (shifts old rows 276..391 down to 277..392,
# Excel automatically extends the sheet dimension/UsedRange accordingly).
$ws.Rows(276).Insert()

# Populate the newly inserted row 276 with the new weekly record.
$ws.Cells.Item(276, 1).Value  = 9
$ws.Cells.Item(276, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(276, 3).Value  = "Metropolitana"
$ws.Cells.Item(276, 4).Value  = 44636
$ws.Cells.Item(276, 5).Value  = 13
$ws.Cells.Item(276, 6).Value  = 100112012
$ws.Cells.Item(276, 7).Value  = "Espinaca"
$ws.Cells.Item(276, 8).Value  = "Sin especificar"
$ws.Cells.Item(276, 9).Value  = "Primera"
$ws.Cells.Item(276, 10).Value = 220
$ws.Cells.Item(276, 11).Value = 10000
$ws.Cells.Item(276, 12).Value = 12000
$ws.Cells.Item(276, 13).Value = 11273
$ws.Cells.Item(276, 14).Value = "`$/cuna 10 kilos"
$ws.Cells.Item(276, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(276, 16).Value = 1127
$ws.Cells.Item(276, 17).Value = 10
$ws.Cells.Item(276, 18).Value = "Hortaliza"
